$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 277
$ws.Range("F3").Value = 472
$ws.Range("F5").Value = 2474
$ws.Range("F9").Value = 1707
$ws.Range("F10").Value = 1707
$ws.Range("F11").Value = 1410
$ws.Range("F13").Value = 1449
$ws.Range("F14").Value = 23
$ws.Range("F15").Value = 33
$ws.Range("F16").Value = 931
$ws.Range("F18").Value = 197
$ws.Range("F20").Value = 7565
$ws.Range("F21").Value = 8598
$ws.Range("F22").Value = 59
$ws.Range("F26").Value = 101
$ws.Range("F27").Value = 274
$ws.Range("F33").Value = 19
$ws.Range("F35").Value = 250
$ws.Range("F36").Value = 29
$ws.Range("F39").Value = 806
$ws.Range("F41").Value = 1375
$ws.Range("F43").Value = 274
$ws.Range("F46").Value = 219
$ws.Range("F47").Value = 5
$ws.Range("F49").Value = 44

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 75
$ws.Range("F6").Value = 23

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 308

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 277
$ws.Range("F4").Value = 472
$ws.Range("F6").Value = 308
$ws.Range("F10").Value = 2474
$ws.Range("F14").Value = 1707
$ws.Range("F15").Value = 1707
$ws.Range("F16").Value = 1410
$ws.Range("F18").Value = 1449
$ws.Range("F19").Value = 23
$ws.Range("F22").Value = 197
$ws.Range("F23").Value = 75
$ws.Range("F25").Value = 7565
$ws.Range("F26").Value = 8598
$ws.Range("F27").Value = 59
$ws.Range("F29").Value = 101
$ws.Range("F30").Value = 274
$ws.Range("F35").Value = 19
$ws.Range("F37").Value = 250
$ws.Range("F38").Value = 29
$ws.Range("F41").Value = 806
$ws.Range("F43").Value = 1375
$ws.Range("F45").Value = 274
$ws.Range("F47").Value = 219
$ws.Range("F50").Value = 44
